$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Force text storage (not boolean/date auto-conversion) to match the
# original author's intent of plain text shared strings.
$ws.Range("B7").NumberFormat = "@"
$ws.Range("B8").NumberFormat = "@"
$ws.Range("B14").NumberFormat = "@"

# Experimental flag flipped from false to true
$ws.Range("B7").Value = "true"

# Date updated
$ws.Range("B8").Value = "2023-02-16T14:43:10-06:00"

# Case Sensitive value now populated with "true"
$ws.Range("B14").Value = "true"
